# Apply scheduled-runner value updates to the Leve profit tables.
# Each sheet has columns: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 30000
$ws.Range("I46").Value = 30000
$ws.Range("K46").Value = 90000
$ws.Range("M46").Value = -89881
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 60000
$ws.Range("N48").Value = -60584
$ws.Range("H56").Value = 20000
$ws.Range("J56").Value = 20000
$ws.Range("L56").Value = 60000
$ws.Range("N56").Value = -61068
$ws.Range("H60").Value = 30000
$ws.Range("I60").Value = 30000
$ws.Range("K60").Value = 90000
$ws.Range("M60").Value = -89516
$ws.Range("H86").Value = 2998
$ws.Range("I86").Value = 2998
$ws.Range("K86").Value = 2998
$ws.Range("M86").Value = -1875
$ws.Range("H89").Value = 2998
$ws.Range("I89").Value = 2998
$ws.Range("K89").Value = 14990
$ws.Range("M89").Value = -9374
$ws.Range("H107").Value = 176.66667
$ws.Range("I107").Value = 192
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 192
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = 1728
$ws.Range("N107").Value = -3940
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H132").Value = 2741.5715
$ws.Range("I132").Value = 2741.5715
$ws.Range("K132").Value = 8224.7145
$ws.Range("M132").Value = -5694.7145
$ws.Range("H138").Value = 1336.375
$ws.Range("J138").Value = 3000
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5774.9
$ws.Range("I32").Value = 5543.625
$ws.Range("K32").Value = 5543.625
$ws.Range("M32").Value = -5256.625
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H88").Value = 336.375
$ws.Range("I88").Value = 336.375
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 336.375
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 69.625
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 336.375
$ws.Range("I91").Value = 336.375
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 336.375
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 1067.625
$ws.Range("N91").ClearContents()
$ws.Range("H97").Value = 1533
$ws.Range("I97").Value = 1533
$ws.Range("K97").Value = 1533
$ws.Range("M97").Value = -1037
$ws.Range("H110").Value = 450
$ws.Range("I110").Value = 450
$ws.Range("K110").Value = 450
$ws.Range("M110").Value = 1595
$ws.Range("H132").Value = 3737.4285
$ws.Range("I132").Value = 3255.6924
$ws.Range("K132").Value = 9767.0772
$ws.Range("M132").Value = -7237.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 486.16666
$ws.Range("I94").Value = 486.16666
$ws.Range("K94").Value = 486.16666
$ws.Range("M94").Value = -35.16665999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1972.5
$ws.Range("I31").Value = 1972.5
$ws.Range("K31").Value = 1972.5
$ws.Range("M31").Value = -1677.5
$ws.Range("H34").Value = 1972.5
$ws.Range("I34").Value = 1972.5
$ws.Range("K34").Value = 1972.5
$ws.Range("M34").Value = -1770.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 282.9
$ws.Range("I23").Value = 42.5
$ws.Range("J23").Value = 443.16666
$ws.Range("K23").Value = 127.5
$ws.Range("L23").Value = 1329.49998
$ws.Range("M23").Value = 107.5
$ws.Range("N23").Value = -1799.49998
$ws.Range("H38").Value = 1917.3334
$ws.Range("I38").Value = 1595.4
$ws.Range("J38").Value = 2319.75
$ws.Range("K38").Value = 4786.200000000001
$ws.Range("L38").Value = 6959.25
$ws.Range("M38").Value = -4439.200000000001
$ws.Range("N38").Value = -7653.25
$ws.Range("H114").Value = 508.66666
$ws.Range("I114").Value = 528
$ws.Range("J114").Value = 499
$ws.Range("K114").Value = 1584
$ws.Range("L114").Value = 1497
$ws.Range("M114").Value = 1670
$ws.Range("N114").Value = -8005
$ws.Range("H121").Value = 1594.8
$ws.Range("I121").Value = 2000
$ws.Range("J121").Value = 1493.5
$ws.Range("K121").Value = 6000
$ws.Range("L121").Value = 4480.5
$ws.Range("M121").Value = -4690
$ws.Range("N121").Value = -7100.5
$ws.Range("H131").Value = 2565.7144
$ws.Range("I131").Value = 3320
$ws.Range("K131").Value = 9960
$ws.Range("M131").Value = -4920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 200
$ws.Range("K2").Value = 200
$ws.Range("M2").Value = -87
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1002
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -5008
$ws.Range("N83").ClearContents()
$ws.Range("H122").Value = 1801.6666
$ws.Range("I122").Value = 1799.5
$ws.Range("K122").Value = 5398.5
$ws.Range("M122").Value = -2948.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9002
$ws.Range("J2").Value = 9002
$ws.Range("L2").Value = 9002
$ws.Range("N2").Value = -9226
$ws.Range("H43").Value = 28800
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 28800
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 28800
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -29186
$ws.Range("H55").Value = 1000
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 1000
$ws.Range("L55").Value = 1000
$ws.Range("M55").Value = -827
$ws.Range("N55").Value = -1346

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1500
$ws.Range("J107").Value = 1500
$ws.Range("L107").Value = 4500
$ws.Range("N107").Value = -8340
$ws.Range("H132").Value = 500
$ws.Range("I132").Value = 500
$ws.Range("K132").Value = 1500
$ws.Range("M132").Value = 1030
$ws.Range("H136").Value = 2410.2942
$ws.Range("I136").Value = 1398.5333
$ws.Range("K136").Value = 4195.5999
$ws.Range("M136").Value = -1645.5999
